$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.104.91'
$ws.Range("E2").Value = '  -1.30%  '
$ws.Range("D3").Value = '2.984.30'
$ws.Range("E3").Value = '  -2.30%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '500.08'
$ws.Range("E5").Value = '  -4.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.52'
$ws.Range("E6").Value = '  -3.49%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -4.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.28'
$ws.Range("E9").Value = '  -5.38%  '
$ws.Range("E10").Value = '  -4.34%  '
$ws.Range("E11").Value = '  -4.12%  '
$ws.Range("D12").Value = '3.494.78'
$ws.Range("E12").Value = '  -2.33%  '
$ws.Range("E13").Value = '  -2.47%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.07'
$ws.Range("E14").Value = '  -3.67%  '
$ws.Range("E15").Value = '  -5.70%  '
$ws.Range("D16").Value = '57.135.43'
$ws.Range("E16").Value = '  -1.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.09'
$ws.Range("E17").Value = '  -2.66%  '
$ws.Range("D18").Value = '2.984.72'
$ws.Range("E18").Value = '  -2.49%  '
$ws.Range("E19").Value = '  -3.45%  '
$ws.Range("E20").Value = '  -3.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.52'
$ws.Range("E21").Value = '  -5.09%  '
$ws.Range("E22").Value = '  +0.00%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.75'
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("E24").Value = '  -2.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.98'
$ws.Range("E25").Value = '  -3.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  -5.27%  '
$ws.Range("E28").Value = '  -8.38%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.62'
$ws.Range("E29").Value = '  -4.27%  '
$ws.Range("E30").Value = '  -3.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.78'
$ws.Range("E31").Value = '  -4.04%  '
$ws.Range("E32").Value = '  -6.03%  '
$ws.Range("E33").Value = '  -4.64%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '154.52'
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("E35").Value = '  -3.26%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.79'
$ws.Range("E36").Value = '  -3.31%  '
$ws.Range("E37").Value = '  -6.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '24.42'
$ws.Range("E38").Value = '  -6.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0664'
$ws.Range("E39").Value = '  -5.44%  '
$ws.Range("E40").Value = '  -0.09%  '
$ws.Range("D41").Value = '3.015.03'
$ws.Range("E41").Value = '  -2.44%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("E43").Value = '  -3.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.644'
$ws.Range("E44").Value = '  -2.81%  '
$ws.Range("D45").Value = '2.187.66'
$ws.Range("E45").Value = '  -5.82%  '
$ws.Range("E46").Value = '  -6.17%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.95'
$ws.Range("E47").Value = '  -1.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.931'
$ws.Range("E48").Value = '  -9.58%  '
$ws.Range("E49").Value = '  -4.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.21'
$ws.Range("E50").Value = '  -4.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.78'
$ws.Range("E51").Value = '  -10.87%  '
